$wb = $excel.ActiveWorkbook

# Insert the new "API_Data" worksheet immediately before "Fields_Data"
$fieldsSheet = $wb.Worksheets.Item("Fields_Data")
$apiSheet = $wb.Worksheets.Add($fieldsSheet)
$apiSheet.Name = "API_Data"

$data = New-Object 'object[,]' 10,11
$data[0,0] = 1
$data[0,1] = '\"Eldon Base for stackable storage shelf'
$data[0,2] = 'platinum\"'
$data[0,3] = 'Muhammed MacIntyre'
$data[0,4] = 3
$data[0,5] = -213.25
$data[0,6] = 38.94
$data[0,7] = 35
$data[0,8] = 'Nunavut'
$data[0,9] = 'Storage & Organization'
$data[0,10] = 0.8
$data[1,0] = 2
$data[1,1] = '\"1.7 Cubic Foot Compact \"\"Cube\"\" Office Refrigerators\"'
$data[1,2] = 'Barry French'
$data[1,3] = 293
$data[1,4] = 457.81
$data[1,5] = 208.16
$data[1,6] = 68.02
$data[1,7] = 'Nunavut'
$data[1,8] = 'Appliances'
$data[1,9] = 0.58
$data[1,10] = $null
$data[2,0] = 3
$data[2,1] = '\"Cardinal Slant-D� Ring Binder'
$data[2,2] = 'Heavy Gauge Vinyl\"'
$data[2,3] = 'Barry French'
$data[2,4] = 293
$data[2,5] = 46.71
$data[2,6] = 8.69
$data[2,7] = 2.99
$data[2,8] = 'Nunavut'
$data[2,9] = 'Binders and Binder Accessories'
$data[2,10] = 0.39
$data[3,0] = 4
$data[3,1] = 'R380'
$data[3,2] = 'Clay Rozendal'
$data[3,3] = 483
$data[3,4] = 1198.97
$data[3,5] = 195.99
$data[3,6] = 3.99
$data[3,7] = 'Nunavut'
$data[3,8] = 'Telephones and Communication'
$data[3,9] = 0.58
$data[3,10] = $null
$data[4,0] = 5
$data[4,1] = 'Holmes HEPA Air Purifier'
$data[4,2] = 'Carlos Soltero'
$data[4,3] = 515
$data[4,4] = 30.94
$data[4,5] = 21.78
$data[4,6] = 5.94
$data[4,7] = 'Nunavut'
$data[4,8] = 'Appliances'
$data[4,9] = 0.5
$data[4,10] = $null
$data[5,0] = 6
$data[5,1] = 'G.E. Longer-Life Indoor Recessed Floodlight Bulbs'
$data[5,2] = 'Carlos Soltero'
$data[5,3] = 515
$data[5,4] = 4.43
$data[5,5] = 6.64
$data[5,6] = 4.95
$data[5,7] = 'Nunavut'
$data[5,8] = 'Office Furnishings'
$data[5,9] = 0.37
$data[5,10] = $null
$data[6,0] = 7
$data[6,1] = '\"Angle-D Binders with Locking Rings'
$data[6,2] = 'Label Holders\"'
$data[6,3] = 'Carl Jackson'
$data[6,4] = 613
$data[6,5] = -54.04
$data[6,6] = 7.3
$data[6,7] = 7.72
$data[6,8] = 'Nunavut'
$data[6,9] = 'Binders and Binder Accessories'
$data[6,10] = 0.38
$data[7,0] = 8
$data[7,1] = '\"SAFCO Mobile Desk Side File'
$data[7,2] = 'Wire Frame\"'
$data[7,3] = 'Carl Jackson'
$data[7,4] = 613
$data[7,5] = 127.7
$data[7,6] = 42.76
$data[7,7] = 6.22
$data[7,8] = 'Nunavut'
$data[7,9] = 'Storage & Organization'
$data[7,10] = $null
$data[8,0] = 9
$data[8,1] = '\"SAFCO Commercial Wire Shelving'
$data[8,2] = 'Black\"'
$data[8,3] = 'Monica Federle'
$data[8,4] = 643
$data[8,5] = -695.26
$data[8,6] = 138.14
$data[8,7] = 35
$data[8,8] = 'Nunavut'
$data[8,9] = 'Storage & Organization'
$data[8,10] = $null
$data[9,0] = 10
$data[9,1] = 'Xerox 198'
$data[9,2] = 'Dorothy Badders'
$data[9,3] = 678
$data[9,4] = -226.36
$data[9,5] = 4.98
$data[9,6] = 8.33
$data[9,7] = 'Nunavut'
$data[9,8] = 'Paper'
$data[9,9] = 0.38
$data[9,10] = $null
$apiSheet.Range("A1:K10").Value = $data
$apiSheet.Range("B10").Select()
